# Rewrites the neuroscience essay into a chemistry essay, per the commit diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Replace-Text: not found: $old"
    }
}

function Replace-WholeWord($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Replace-WholeWord: not found: $old"
    }
}

# Returns the document-character offset just after the (unique) search text.
function Get-EndOf($searchText) {
    $r = $d.Content
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Get-EndOf: not found: $searchText"
    }
    return $r.End
}

# Apply the body's standard Calibri / black / $size run formatting to [startPos, endPos).
# $size of 0 means "leave the (inherited) default size alone" - matches the Summary
# paragraph's runs, which carry no explicit <w:sz>.
function Format-Range($startPos, $endPos, $size) {
    $fr = $d.Range($startPos, $endPos)
    $fr.Font.Name = "Calibri"
    if ($size -gt 0) {
        $fr.Font.Size = $size
    }
    $fr.Font.Color = 0
}

# Insert a manual line break (<w:br/>) at pos, formatted, returning the new position.
function Insert-LineBreak($pos, $size) {
    $ins = $d.Range($pos, $pos)
    $ins.InsertBreak(6)
    $endPos = $pos + 1
    Format-Range $pos $endPos $size
    return $endPos
}

# Insert formatted text at pos, returning the new position.
function Insert-RunAfter($pos, $text, $size) {
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($text)
    $endPos = $pos + $text.Length
    Format-Range $pos $endPos $size
    return $endPos
}

# ---------------------------------------------------------------------------
# Title / byline / e-mail
# ---------------------------------------------------------------------------

Replace-Text "Unraveling the Mysteries of the Mind: Neuroscience's Journey" "The Enchanting Realm of Chemistry: Unveiling the Secrets of Matter"
Replace-Text " Isabella Armstrong" " Clara Thompson"
Replace-Text "iarmstrong@sciencereview" "clarathompson0304@institute"
Replace-WholeWord "org" "edu"

# ---------------------------------------------------------------------------
# Body paragraph (introduction)
# ---------------------------------------------------------------------------

Replace-Text "In the vast expanse of human knowledge, the human brain stands as one of the most captivating enigmas" "Chemistry, a compelling branch of natural science, unveils the hidden secrets of the substances that constitute our universe"
Replace-Text " This intricate organ, composed of billions of neurons and trillions of connections, orchestrates our thoughts, emotions, and actions" " This captivating field of study delves into the composition, structure, and transformation of matter"
Replace-Text " Over the centuries, philosophers and scientists have sought to unravel the mysteries of the mind, peeling back layer upon layer to reveal its hidden depths" " It explores how elements and compounds interact, unraveling the intricate dance of atoms and molecules"
Replace-Text " In recent decades, the field of neuroscience has emerged as a beacon of discovery, shedding light on the intricate workings of the brain and its profound influence on our lives" " Embarking on a chemical journey, high school students can discover the mesmerizing world of elements, compounds, and their intricate interactions"

Replace-Text "From the pioneering research of Santiago Ramon y Cajal, who revealed the intricate architecture of neurons, to the groundbreaking insights of Eric Kandel, who unraveled the molecular mechanisms of memory formation, neuroscience has made remarkable strides in unraveling the complexities of the brain" "From the smallest subatomic particles to the vast expanse of the cosmos, chemistry orchestrates a symphony of elemental harmony"
Replace-Text " Modern advancements in neuroimaging techniques, such as functional magnetic resonance imaging (fMRI) and electroencephalography (EEG), have enabled scientists to visualize brain activity in real time, revealing patterns of neural communication underlying our cognitive functions and emotional experiences" " It provides the foundation for comprehending the properties of substances, the mechanisms of chemical reactions, and the dynamics of energy transfer"

# New: ". By examining chemical processes..." inserted right after "...energy transfer" and
# before the "." that used to directly follow "...emotional experiences".
$p = Get-EndOf "the dynamics of energy transfer"
$p = Insert-RunAfter $p "." 12
$p = Insert-RunAfter $p " By examining chemical processes, students gain a deeper appreciation for the interplay between structure and function, revealing the elegance and interconnectedness of the natural world" 12

Replace-Text "Furthermore, the advent of optogenetics, a technique that allows for the precise control of neuronal activity using light, has opened up new avenues for investigating the causal relationships between neural circuits and behavior" "Introduction Continued:"

# New: a line break + "Chemistry shapes civilization's progress..." sentence.
$p = Get-EndOf "Introduction Continued:"
$p = Insert-LineBreak $p 12
$p = Insert-RunAfter $p "Chemistry shapes civilization's progress, underpinning countless industries ranging from medicine to materials science" 12

Replace-Text " These groundbreaking discoveries have laid the foundation for a deeper understanding of the brain, revolutionizing our perspectives on consciousness, decision-making, and mental health disorders" " It enables the development of life-saving drugs, durable materials, and sustainable energy sources"

# New block of sentences appended at the end of the introduction paragraph.
$p = Get-EndOf "durable materials, and sustainable energy sources"
$p = Insert-RunAfter $p "." 12
$p = Insert-RunAfter $p " Delving into chemistry enriches our understanding of the world around us, bridging the gap between the abstract and the tangible" 12
$p = Insert-RunAfter $p "." 12
$p = Insert-RunAfter $p " It nurtures a spirit of inquiry, fostering problem-solving skills and critical thinking abilities that extend far beyond the classroom, empowering students to analyze complex phenomena and make informed decisions throughout their lives" 12
$p = Insert-RunAfter $p "." 12
$p = Insert-LineBreak $p 12
$p = Insert-LineBreak $p 12
$p = Insert-RunAfter $p "Introduction Concluded:" 12
$p = Insert-LineBreak $p 12
$p = Insert-RunAfter $p "The study of chemistry ignites a sense of wonder and curiosity, encouraging students to contemplate the nature of matter and its intricate workings" 12
$p = Insert-RunAfter $p "." 12
$p = Insert-RunAfter $p " It empowers them with the knowledge and skills necessary to navigate the complexities of the modern world, addressing challenges such as climate change, resource scarcity, and disease" 12
$p = Insert-RunAfter $p "." 12
$p = Insert-RunAfter $p " Chemistry " 12
$p = Insert-RunAfter $p "inspires creativity and innovation, cultivating a generation of problem-solvers who possess the ability to shape a sustainable and equitable future for all" 12

# ---------------------------------------------------------------------------
# Summary paragraph
# ---------------------------------------------------------------------------

Replace-Text "The field of neuroscience has illuminated the intricate workings of the human brain, revealing the neural underpinnings of our thoughts, emotions, and behaviors" "Chemistry offers a captivating journey into the realm of matter, revealing the secrets of its composition, structure, and transformation"
Replace-Text " Through the dedication of pioneering researchers and the advent of sophisticated neurotechnologies, we have gained unprecedented insights into the complexities of the mind" " It enriches our understanding of the natural world, fostering critical thinking skills and inspiring creativity"

# Collapse the four closing runs (incl. the old page-break split) into the new single sentence.
$r = $d.Content
$startOk = $r.Find.Execute(" This knowledge holds the potential to transform our understanding of mental illness, optimize cognitive performance, and enhance ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $startOk) { throw "summary tail start not found" }
$startPos = $r.Start

$r2 = $d.Content
$endOk = $r2.Find.Execute("As neuroscience continues to forge ahead, we can anticipate even more remarkable discoveries that will deepen our appreciation for the marvel that is the human brain", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $endOk) { throw "summary tail end not found" }
$endPos = $r2.End

$whole = $d.Range($startPos, $endPos)
$newText = " By unraveling the intricacies of chemical interactions, students gain the knowledge and skills they need to address global challenges, shaping a sustainable and equitable future"
$whole.Text = $newText
$newEnd = $startPos + $newText.Length
Format-Range $startPos $newEnd 0

# ---------------------------------------------------------------------------
# Trailing empty paragraph before the section break
# ---------------------------------------------------------------------------

$endOfDoc = $d.Content.End
$tail = $d.Range($endOfDoc, $endOfDoc)
$tail.InsertParagraphAfter()

Write-Output "done"
